# Applies the "details.xlsx" row-2 data update described in the commit.
# The sheet holds a single header row (row 1) plus one data row (row 2);
# this updates that data row's values, blanks a few fields out, shifts the
# last two populated columns one slot to the right, and widens the sheet's
# declared used-range to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Personal / contact info -------------------------------------------------
$ws.Range("A2").Value = 23
$ws.Range("B2").Value = "Jaligama"
$ws.Range("C2").Value = "Prabhu"
$ws.Range("D2").Value = "Jaligama"
# E2 (phone number) must stay text-typed even though it looks numeric.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "8155933548"
$ws.Range("E2").ClearFormats()
# F2 (email) and G2 (sex) are unchanged.

# --- Age / DOB are cleared out ----------------------------------------------
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

# --- Address / State ---------------------------------------------------------
$ws.Range("J2").Value = "5035 s east end S2402`ns2402"
$ws.Range("K2").Value = "ILLINOIS"
# L2 (Risk Class) is unchanged.

# --- Face amount cleared -------------------------------------------------
$ws.Range("M2").ClearContents()
# N2, O2, P2, Q2 unchanged.

# --- Maximum Monthly Benefit cleared -----------------------------------------
$ws.Range("R2").ClearContents()
# S2, T2, U2 unchanged.

# --- Benefit Durations / Inflation Benefit Option shift one column right ----
$ws.Range("V2").ClearContents()
$ws.Range("W2").Value = "2 Years"
$ws.Range("X2").Value = "None"

# --- Widen the sheet's used range to A1:AA2, matching the source file ------
$ws.Range("AA2").Interior.Pattern = -4142
